# Rik Activity 2019.xlsx -- "Input by date instead of range. Updated data thru morning Apr 3"
#
# Appends four new activity log rows (165-168) to the "2019" sheet / Table2,
# adds the new "Pasta + cheese + bread" comment string, and moves the active
# selection to the next empty row (A169), mirroring how the workbook owner
# continued typing new entries directly below the existing table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2019")
$tbl = $ws.ListObjects.Item("Table2")

$templateRow = 164
$xlPasteFormats = -4122

function Add-ActivityRow {
    param(
        [int]$RowNum,
        [double]$Start,
        $End,
        [string]$Activity,
        $Comment
    )

    # Grow the table by one row so its range / autofilter / sheet dimension
    # all get extended the way Excel does when you type into the row right
    # below a table.
    $null = $tbl.ListRows.Add()

    # Copy the Start/End/Activity/Comment formatting down from the last
    # existing data row (E is left alone - it gets its own formula below).
    $ws.Range("A${templateRow}:D${templateRow}").Copy()
    $ws.Range("A${RowNum}:D${RowNum}").PasteSpecial($xlPasteFormats)
    $excel.CutCopyMode = 0

    $ws.Range("A$RowNum").Value = $Start

    if ($End -ne $null) {
        $ws.Range("A${templateRow}").Copy()
        $ws.Range("B$RowNum").PasteSpecial($xlPasteFormats)
        $excel.CutCopyMode = 0
        $ws.Range("B$RowNum").Value = $End
    } else {
        $ws.Range("B$RowNum").Clear()
    }

    $ws.Range("C$RowNum").Value = $Activity

    if ($Comment -ne $null) {
        $ws.Range("D$RowNum").Value = $Comment
    } else {
        $ws.Range("D$RowNum").Clear()
    }

    $ws.Range("E$RowNum").Formula = '=IF(Table2[[#This Row],[Activity]]="Sleep",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,"NA")'
}

Add-ActivityRow 165 43557.886805555558 $null               "Food"  "Pasta + cheese + bread"
Add-ActivityRow 166 43557.928136574075 43558.270138888889  "Sleep" $null
Add-ActivityRow 167 43557.718055555553 $null               "Food"  "Nuts"
Add-ActivityRow 168 43558.28125        $null               "Food"  "Latte"

# Matches the author's last selection after typing the new rows.
$ws.Range("A169").Select()
